$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New/updated "close" values for B2:B133 (row index -> value)
$values = @(
    1.19685, 1.16442, 1.15956, 1.17902, 1.19685, 1.16118, 1.16767, 1.15956, 1.16605, 1.17577,
    1.17902, 1.18713, 1.20658, 1.22604, 1.22928, 1.21793, 1.21631, 1.22928, 1.25847, 1.23253,
    1.24874, 1.24874, 1.21307, 1.2228, 1.22442, 1.23091, 1.2228, 1.2455, 1.26496, 1.26009,
    1.25523, 1.24712, 1.24388, 1.25036, 1.26496, 1.2682, 1.26658, 1.25523, 1.22442, 1.22928,
    1.22442, 1.23577, 1.24064, 1.24064, 1.23901, 1.25036, 1.25199, 1.23901, 1.22928, 1.23091,
    1.24388, 1.24388, 1.23901, 1.24064, 1.30874, 1.34279, 1.36712, 1.37684, 1.38982, 1.37522,
    1.39306, 1.39144, 1.38009, 1.40603, 1.38495, 1.3736, 1.44008, 1.43846, 1.44819, 1.4563,
    1.48549, 1.49522, 1.48387, 1.45792, 1.49359, 1.47251, 1.47251, 1.50332, 1.58602, 1.58602,
    1.56981, 1.52116, 1.50495, 1.54548, 1.57305, 1.55197, 1.53738, 1.55359, 1.55683, 1.58764,
    1.57791, 1.58278, 1.58602, 1.55521, 1.53738, 1.53251, 1.51467, 1.52116, 1.53413, 1.55197,
    1.5617, 1.49846, 1.52116, 1.5017, 1.4563, 1.47738, 1.50981, 1.63629, 1.68656, 1.68656,
    1.75304, 1.73034, 1.7579, 1.8179, 1.81304, 1.80169, 1.81304, 1.79194, 1.80493, 1.8633,
    1.88438, 1.86979, 1.84547, 1.89411, 1.92654, 1.91844, 1.82763, 1.84222, 1.91033, 1.88925,
    1.8779, 1.87141
)

# Write column B values starting at row 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Extend column A (index sequence) and copy formatting for the six new rows (128-133)
$lastOldRow = 127
$newLastRow = 133
$srcA = $ws.Range("A$lastOldRow")
$dstA = $ws.Range("A" + ($lastOldRow + 1) + ":A" + $newLastRow)
$srcA.Copy()
$dstA.PasteSpecial(-4122)

for ($row = $lastOldRow + 1; $row -le $newLastRow; $row++) {
    $ws.Cells.Item($row, 1).Value = $row - 2
}

$excel.CutCopyMode = $false

Write-Output "Updated close_validation_scaler values; dimension now A1:B$newLastRow"